# Auto-generated edit script for TC01_Canine_Filter_FileAssoc-diagnosis.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Long Cypher query text blocks (literal here-strings; no $ or ` expansion) ----
$qCasesNew = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
WHERE labels(parent)[0] IN ["diagnosis"] 
MATCH (f)-[*]->(c:case)<--(demo:demographic)
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
 MATCH (samp:sample)-->(c) 
 MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co,demo.patient_age_at_enrollment AS age
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
    coalesce(CASE age % 1 WHEN 0 THEN apoc.convert.toInteger(age) ELSE age END, '') AS Age,
       coalesce(demo.sex, '') AS Sex,
       coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
coalesce(CASE weight % 1 WHEN 0 THEN apoc.convert.toInteger(weight) ELSE weight END, '') AS `Weight (kg)`,
       coalesce(diag.best_response, '') AS `Response to Treatment`,
       coalesce(co.cohort_description, '') AS `Cohort`
'@
$qSamplesUnchanged = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
 MATCH (samp:sample)-->(c) 
 WHERE labels(parent)[0] IN ["diagnosis"]  
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) 
WITH DISTINCT samp AS samp, c, demo, diag
RETURN  coalesce(samp.sample_id, '') AS `Sample ID`, 
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(demo.breed,'') AS Breed , 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(samp.sample_site, '') AS `Sample Site`,
        coalesce(samp.summarized_sample_type, '') AS `Sample Type`,
        coalesce(samp.specific_sample_pathology, '') AS `Pathology/Morphology`,
        coalesce(samp.tumor_grade, '') AS `Tumor Grade`,
        coalesce(samp.sample_chronology, '') AS `Sample Chronology`,
        coalesce(samp.percentage_tumor, '') AS `Percentage Tumor`,
        coalesce(samp.necropsy_sample, '') AS `Necropsy Sample`,
        coalesce(samp.sample_preservation, '') AS `Sample Preservation`
'@
$qFilesNew = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f)-[*]->(samp:sample)
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE labels(parent)[0] IN ["diagnosis"]  
OPTIONAL MATCH (s:study)<--(c)<--(diag:diagnosis)<-[*]-(samp)
WITH
        f, parent, c, demo, diag, s, samp,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent, c, demo, diag, s, samp,
        f.file_size /(1024^i) AS value, 
        10^precision AS factor,
        units[i] as unit
WITH    
        f, parent, c, demo, diag, s, samp, unit,
        round(factor * value)/factor AS size
RETURN 
        coalesce(f.file_name, '') AS `File Name`,
        coalesce(f.file_type, '') AS `File Type`,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
        coalesce(samp.sample_id, '') AS `Sample ID`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(demo.breed,'') AS Breed ,
        coalesce(diag.disease_term,'') AS Diagnosis
'@
$qStudyFilesNew = @'
MATCH (f:file)-->(parent)
MATCH (f:file)-->(s:study)
WITH DISTINCT f, parent
MATCH (s)<--(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)
WHERE labels(parent)[0] IN ["diagnosis"]  
WITH DISTINCT f,  s, c, demo, diag
WITH
        f, c, demo, diag, s,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH    
        f, c, demo, diag, s,
        f.file_size /(1024^i) AS value, 10^precision AS factor,
        units[i] as unit
        WITH    
        f,  c, demo, diag, s, unit,
        round(factor * value)/factor AS size
RETURN DISTINCT
  coalesce(f.file_name, '') AS `File Name`,
  coalesce(f.file_type, '') AS `File Type`,
  coalesce("study", '') AS `Association`,
  coalesce(f.file_description, '') AS `Description`,
  coalesce(f.file_format, '') AS  Format,
  CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
  coalesce(s.clinical_study_designation,'') AS `Study Code`
'@
$qStat = @'
OPTIONAL MATCH (sf:file)-->(ss:study)
WHERE head(labels(ss)) IN ["diagnosis"]
WITH count(distinct sf) AS study_files
OPTIONAL MATCH (f:file)-->(parent)
WHERE head(labels(parent)) IN ["diagnosis"]
OPTIONAL MATCH (f)-[*]->(c:case)
OPTIONAL MATCH (demo:demographic)-->(c)-->(s:study)-->(p:program)
OPTIONAL MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
RETURN
	count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    study_files AS `Study Files`
    
'@

# ---- Row 2: CasesTab (query replaced with age/weight-aware version) ----
$ws.Range("A2").Value = "CasesTab"
$ws.Range("B2").Value = $qCasesNew
$ws.Range("C2").Value = $qStat
$ws.Range("D2").Value = "TC01_Canine_Filter_FileAssoc-diagnosis_Neo4jData.xlsx"
$ws.Range("E2").Value = "TC01_Canine_Filter_FileAssoc-diagnosis_WebData.xlsx"

# ---- Row 3: SamplesTab (query text unchanged, but Column C now uses the StatQuery text) ----
$ws.Range("A3").Value = "SamplesTab"
$ws.Range("B3").Value = $qSamplesUnchanged
$ws.Range("C3").Value = $qStat
$ws.Range("D3").Value = "TC01_Canine_Filter_FileAssoc-diagnosis_Neo4jData.xlsx"
$ws.Range("E3").Value = "TC01_Canine_Filter_FileAssoc-diagnosis_WebData.xlsx"

# ---- Row 4: FilesTab (query replaced with File Type / Sample ID aware version) ----
$ws.Range("A4").Value = "FilesTab"
$ws.Range("B4").Value = $qFilesNew
$ws.Range("C4").Value = $qStat
$ws.Range("D4").Value = "TC01_Canine_Filter_FileAssoc-diagnosis_Neo4jData.xlsx"
$ws.Range("E4").Value = "TC01_Canine_Filter_FileAssoc-diagnosis_WebData.xlsx"

# ---- Row 5 (NEW): StudyFilesTab ----
$ws.Range("A5").Value = "StudyFilesTab"
$ws.Range("B5").Value = $qStudyFilesNew
$ws.Range("C5").Value = $qStat
$ws.Range("D5").Value = "TC01_Canine_Filter_FileAssoc-diagnosis_Neo4jData.xlsx"
$ws.Range("E5").Value = "TC01_Canine_Filter_FileAssoc-diagnosis_WebData.xlsx"

# ---- Row heights ----
$ws.Rows.Item(1).RowHeight = 46.5
$ws.Range("A2:E5").RowHeight = 99.95

# ---- Column widths ----
$ws.Columns.Item(2).ColumnWidth = 67.28515625

# ---- Wrap text for the long query columns (B and C) on data rows ----
$ws.Range("B2:C5").WrapText = $true

# ---- Selection / view ----
$ws.Range("C2").Select()
